$d = $word.ActiveDocument

# --- Title: "Factor Quota values" -> "Factor ABOVE TAC values" ---
$d.Content.Find.Execute("Table 19. Factor Quota values per country.", $true, $false, $false, $false, $false, $true, 0, $false, "Table 19. Factor ABOVE TAC values per country.", 1) | Out-Null

# --- Table header labels: QUOTAS -> ABOVE.TAC (both columns) ---
$t = $d.Tables.Item(1)
$h2 = $t.Cell(1, 2); $hr2 = $d.Range($h2.Range.Start, $h2.Range.End)
$hr2.Find.Execute("QUOTAS", $true, $false, $false, $false, $false, $true, 0, $false, "ABOVE.TAC", 1) | Out-Null
$h3 = $t.Cell(1, 3); $hr3 = $d.Range($h3.Range.Start, $h3.Range.End)
$hr3.Find.Execute("QUOTAS", $true, $false, $false, $false, $false, $true, 0, $false, "ABOVE.TAC", 1) | Out-Null

# --- Data values per country row (column 2 = cod, column 3 = hake) ---
# BE
$c2 = $t.Cell(2, 2); $rng2 = $d.Range($c2.Range.Start, $c2.Range.End)
$rng2.Find.Execute("0.018", $true, $false, $false, $false, $false, $true, 0, $false, "0.958", 1) | Out-Null
$c3 = $t.Cell(2, 3); $rng3 = $d.Range($c3.Range.Start, $c3.Range.End)
$rng3.Find.Execute("0.011", $true, $false, $false, $false, $false, $true, 0, $false, "0.976", 1) | Out-Null
# DK
$c2 = $t.Cell(3, 2); $rng2 = $d.Range($c2.Range.Start, $c2.Range.End)
$rng2.Find.Execute("0.350", $true, $false, $false, $false, $false, $true, 0, $false, "0.000", 1) | Out-Null
$c3 = $t.Cell(3, 3); $rng3 = $d.Range($c3.Range.Start, $c3.Range.End)
$rng3.Find.Execute("0.091", $true, $false, $false, $false, $false, $true, 0, $false, "0.967", 1) | Out-Null
# DE
$c2 = $t.Cell(4, 2); $rng2 = $d.Range($c2.Range.Start, $c2.Range.End)
$rng2.Find.Execute("0.303", $true, $false, $false, $false, $false, $true, 0, $false, "0.739", 1) | Out-Null
$c3 = $t.Cell(4, 3); $rng3 = $d.Range($c3.Range.Start, $c3.Range.End)
$rng3.Find.Execute("0.004", $true, $false, $false, $false, $false, $true, 0, $false, "0.976", 1) | Out-Null
# EE
$c2 = $t.Cell(5, 2); $rng2 = $d.Range($c2.Range.Start, $c2.Range.End)
$rng2.Find.Execute("0.024", $true, $false, $false, $false, $false, $true, 0, $false, "0.936", 1) | Out-Null
$c3 = $t.Cell(5, 3); $rng3 = $d.Range($c3.Range.Start, $c3.Range.End)
$rng3.Find.Execute("0.000", $true, $false, $false, $false, $false, $true, 0, $false, "0.977", 1) | Out-Null
# IE
$c2 = $t.Cell(6, 2); $rng2 = $d.Range($c2.Range.Start, $c2.Range.End)
$rng2.Find.Execute("0.003", $true, $false, $false, $false, $false, $true, 0, $false, "0.917", 1) | Out-Null
$c3 = $t.Cell(6, 3); $rng3 = $d.Range($c3.Range.Start, $c3.Range.End)
$rng3.Find.Execute("0.059", $true, $false, $false, $false, $false, $true, 0, $false, "0.970", 1) | Out-Null
# ES
$c2 = $t.Cell(7, 2); $rng2 = $d.Range($c2.Range.Start, $c2.Range.End)
$rng2.Find.Execute("0.277", $true, $false, $false, $false, $false, $true, 0, $false, "1.000", 1) | Out-Null
$c3 = $t.Cell(7, 3); $rng3 = $d.Range($c3.Range.Start, $c3.Range.End)
$rng3.Find.Execute("0.716", $true, $false, $false, $false, $false, $true, 0, $false, "0.844", 1) | Out-Null
# FR
$c2 = $t.Cell(8, 2); $rng2 = $d.Range($c2.Range.Start, $c2.Range.End)
$rng2.Find.Execute("0.091", $true, $false, $false, $false, $false, $true, 0, $false, "0.861", 1) | Out-Null
$c3 = $t.Cell(8, 3); $rng3 = $d.Range($c3.Range.Start, $c3.Range.End)
$rng3.Find.Execute("1.000", $true, $false, $false, $false, $false, $true, 0, $false, "0.930", 1) | Out-Null
# LV
$c2 = $t.Cell(9, 2); $rng2 = $d.Range($c2.Range.Start, $c2.Range.End)
$rng2.Find.Execute("0.092", $true, $false, $false, $false, $false, $true, 0, $false, "0.824", 1) | Out-Null
$c3 = $t.Cell(9, 3); $rng3 = $d.Range($c3.Range.Start, $c3.Range.End)
$rng3.Find.Execute("0.000", $true, $false, $false, $false, $false, $true, 0, $false, "0.977", 1) | Out-Null
# LT
$c2 = $t.Cell(10, 2); $rng2 = $d.Range($c2.Range.Start, $c2.Range.End)
$rng2.Find.Execute("0.060", $true, $false, $false, $false, $false, $true, 0, $false, "0.874", 1) | Out-Null
$c3 = $t.Cell(10, 3); $rng3 = $d.Range($c3.Range.Start, $c3.Range.End)
$rng3.Find.Execute("0.000", $true, $false, $false, $false, $false, $true, 0, $false, "0.977", 1) | Out-Null
# NL
$c2 = $t.Cell(11, 2); $rng2 = $d.Range($c2.Range.Start, $c2.Range.End)
$rng2.Find.Execute("0.058", $true, $false, $false, $false, $false, $true, 0, $false, "0.912", 1) | Out-Null
$c3 = $t.Cell(11, 3); $rng3 = $d.Range($c3.Range.Start, $c3.Range.End)
$rng3.Find.Execute("0.009", $true, $false, $false, $false, $false, $true, 0, $false, "0.976", 1) | Out-Null
# PL
$c2 = $t.Cell(12, 2); $rng2 = $d.Range($c2.Range.Start, $c2.Range.End)
$rng2.Find.Execute("0.341", $true, $false, $false, $false, $false, $true, 0, $false, "0.664", 1) | Out-Null
$c3 = $t.Cell(12, 3); $rng3 = $d.Range($c3.Range.Start, $c3.Range.End)
$rng3.Find.Execute("0.000", $true, $false, $false, $false, $false, $true, 0, $false, "0.977", 1) | Out-Null
# PT
$c2 = $t.Cell(13, 2); $rng2 = $d.Range($c2.Range.Start, $c2.Range.End)
$rng2.Find.Execute("0.056", $true, $false, $false, $false, $false, $true, 0, $false, "0.984", 1) | Out-Null
$c3 = $t.Cell(13, 3); $rng3 = $d.Range($c3.Range.Start, $c3.Range.End)
$rng3.Find.Execute("0.086", $true, $false, $false, $false, $false, $true, 0, $false, "0.979", 1) | Out-Null
# FI
$c2 = $t.Cell(14, 2); $rng2 = $d.Range($c2.Range.Start, $c2.Range.End)
$rng2.Find.Execute("0.019", $true, $false, $false, $false, $false, $true, 0, $false, "0.901", 1) | Out-Null
$c3 = $t.Cell(14, 3); $rng3 = $d.Range($c3.Range.Start, $c3.Range.End)
$rng3.Find.Execute("0.000", $true, $false, $false, $false, $false, $true, 0, $false, "0.977", 1) | Out-Null
# SE
$c2 = $t.Cell(15, 2); $rng2 = $d.Range($c2.Range.Start, $c2.Range.End)
$rng2.Find.Execute("0.250", $true, $false, $false, $false, $false, $true, 0, $false, "0.479", 1) | Out-Null
$c3 = $t.Cell(15, 3); $rng3 = $d.Range($c3.Range.Start, $c3.Range.End)
$rng3.Find.Execute("0.004", $true, $false, $false, $false, $false, $true, 0, $false, "0.976", 1) | Out-Null
